# "Generate Report for Handback" - refresh the localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text lives in one shared string used by Overview!E2/F2 and by the
#     "Status" column on each language sheet, so setting it once per sheet
#     updates every cell that shares it)
#   - The per-language "Latest Handback DateTime" is refreshed to the time the
#     handback report was generated
#   - The stale "handback file is not the latest" Error Detail is cleared now
#     that the handback is in sync
#   - A couple of columns are widened so the longer status/blank text isn't
#     truncated

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Handed back: in sync with en-US"
$ws.Range("F2").Value = "Handed back: in sync with en-US"
$ws.Columns.Item(5).ColumnWidth = 29.15
$ws.Columns.Item(6).ColumnWidth = 29.15

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("K2").Value = "2016-08-30 15:04:44"
$ws.Range("P2").Value = ""
$ws.Columns.Item(3).ColumnWidth = 29.15
$ws.Columns.Item(16).ColumnWidth = 12.83

# ---- de-de sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("K2").Value = "2016-08-30 15:04:51"
$ws.Range("P2").Value = ""
$ws.Columns.Item(3).ColumnWidth = 29.15
$ws.Columns.Item(16).ColumnWidth = 12.83
